$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = "cloumn1"
$ws.Range("C3").Value = "column2"
$ws.Range("D3").Value = "cloumn3"
$ws.Range("E3").Value = "column4"

$ws.Range("E3").Select() | Out-Null
